{"js": "// Update the worked division answers in the single table of the\n// document. The table has 20 rows x 5 columns; only rows 0, 4, 8, 12,\n// 16 contain text (the others are spacer rows). We overwrite each\n// populated cell's text with its new value, in table (row, col) order,\n// matching the order of replacements in the source diff. Cell\n// formatting (font / size / paragraph alignment) is left untouched\n// because we only change the cell's text value.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New text for each populated row, left-to-right by column.\nconst newRowValues = {\n  0: [\"27\u00f73=9, 0\", \"90\u00f72=45, 0\", \"18\u00f73=6, 0\", \"14\u00f78=1, 6\", \"48\u00f76=8, 0\"],\n  4: [\"13\u00f75=2, 3\", \"40\u00f73=13, 1\", \"20\u00f79=2, 2\", \"97\u00f75=19, 2\", \"16\u00f78=2, 0\"],\n  8: [\"71\u00f72=35, 1\", \"39\u00f76=6, 3\", \"54\u00f79=6, 0\", \"49\u00f78=6, 1\", \"71\u00f72=35, 1\"],\n  12: [\"78\u00f79=8, 6\", \"49\u00f72=24, 1\", \"99\u00f77=14, 1\", \"91\u00f79=10, 1\", \"93\u00f77=13, 2\"],\n  16: [\"32\u00f75=6, 2\", \"35\u00f73=11, 2\", \"35\u00f75=7, 0\", \"21\u00f76=3, 3\", \"41\u00f76=6, 5\"],\n};\n\nfor (const rowIndexStr of Object.keys(newRowValues)) {\n  const rowIndex = Number(rowIndexStr);\n  const values = newRowValues[rowIndex];\n  for (let colIndex = 0; colIndex < values.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    cell.value = values[colIndex];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worked division answers in the single table of the\n# document. The table has 20 rows x 5 columns (1-based in the Word\n# object model); only rows 1, 5, 9, 13, 17 contain text (the others are\n# spacer rows). We overwrite each populated cell's Range.Text with its\n# new value, leaving the existing run formatting (font / size /\n# paragraph alignment) untouched since we only change the text content\n# of the cell's range.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowsData = @{\n    1  = @(\"27\u00f73=9, 0\", \"90\u00f72=45, 0\", \"18\u00f73=6, 0\", \"14\u00f78=1, 6\", \"48\u00f76=8, 0\")\n    5  = @(\"13\u00f75=2, 3\", \"40\u00f73=13, 1\", \"20\u00f79=2, 2\", \"97\u00f75=19, 2\", \"16\u00f78=2, 0\")\n    9  = @(\"71\u00f72=35, 1\", \"39\u00f76=6, 3\", \"54\u00f79=6, 0\", \"49\u00f78=6, 1\", \"71\u00f72=35, 1\")\n    13 = @(\"78\u00f79=8, 6\", \"49\u00f72=24, 1\", \"99\u00f77=14, 1\", \"91\u00f79=10, 1\", \"93\u00f77=13, 2\")\n    17 = @(\"32\u00f75=6, 2\", \"35\u00f73=11, 2\", \"35\u00f75=7, 0\", \"21\u00f76=3, 3\", \"41\u00f76=6, 5\")\n}\n\nforeach ($rowIndex in $rowsData.Keys) {\n    $values = $rowsData[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
